# Revert "20/10/21 Fixed Build Issues"
#
# On slide 3 ("Prayer Requests (Streamelements)"):
#  - remove the redundant "Prayer Requests- 代禱" text box (id 7)
#  - remove the "請 關 您 的 手 提 電 話" text box (id 9)
#  - remove the "PLEASE SWITCH OFF YOUR MOBILE PHONE" text box (id 10)
#  - remove the "no mobile phones" sign picture (id 11, rId2 -> image1.jpeg)
#  - make the full-slide background picture (id 1026) point at the same
#    image that the removed sign used (rId2 -> image1.jpeg) instead of
#    rId3 -> image2.jpeg

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

function Get-ShapeById($slide, $id) {
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $sh = $slide.Shapes.Item($j)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# The small "no mobile phones" picture already embeds the image
# (rId2 / image1.jpeg) that the big background picture needs to switch
# to. Duplicate it first so we keep a live reference to that image,
# then stretch/crop the duplicate to cover the full slide exactly like
# the background picture did.
$signPic = Get-ShapeById $s 11
$dup = $signPic.Duplicate()
$newPic = $dup.Item(1)

$newPic.PictureFormat.CropLeft = 0
$newPic.PictureFormat.CropTop = 0
$newPic.PictureFormat.CropRight = 0
$newPic.PictureFormat.CropBottom = 0

$newPic.Left = 0
$newPic.Top = 0
$newPic.Width = 960
$newPic.Height = 540

$newPic.Name = "Picture 2"
$newPic.AlternativeText = ""

# Remove the old full-slide background picture (rId3 / image2.jpeg).
$oldBackground = Get-ShapeById $s 1026
$oldBackground.Delete()

# Remove the original small "no mobile phones" picture and its
# accompanying text boxes.
$signPic.Delete()

$textBox7 = Get-ShapeById $s 7
if ($textBox7 -ne $null) { $textBox7.Delete() }

$textBox9 = Get-ShapeById $s 9
if ($textBox9 -ne $null) { $textBox9.Delete() }

$textBox10 = Get-ShapeById $s 10
if ($textBox10 -ne $null) { $textBox10.Delete() }
